{"js": "// Merge the word-by-word runs of the Title, Author and Abstract\n// paragraphs into a single run each (same visible text, just\n// collapsing the run structure that split every word/space into\n// its own <w:r>).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst targets = {\n  \"Title\": \"Answers: Using the quadratic formula\",\n  \"Author\": \"Tom Coleman\",\n  \"Abstract\": \"Answers to questions relating to the guide on using the quadratic formula.\"\n};\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const style = para.style;\n  if (Object.prototype.hasOwnProperty.call(targets, style)) {\n    const expected = targets[style];\n    // Re-write the paragraph's whole range with the same text so the\n    // many single-word runs collapse into one run.\n    para.getRange(\"Whole\").insertText(expected, \"Replace\");\n    delete targets[style];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Merge the word-by-word runs of the Title, Author and Abstract\n# paragraphs into a single run each (same visible text, just\n# collapsing the run structure that split every word/space into\n# its own run).\n$d = $word.ActiveDocument\n\n$targets = @{\n  \"Title\"    = \"Answers: Using the quadratic formula\"\n  \"Author\"   = \"Tom Coleman\"\n  \"Abstract\" = \"Answers to questions relating to the guide on using the quadratic formula.\"\n}\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $styleName = $p.Style.NameLocal\n  if ($targets.ContainsKey($styleName)) {\n    $target = $targets[$styleName]\n\n    # Setting Range.Text to the exact same text it already holds is a\n    # no-op for the engine (no real text change => runs are left\n    # alone), so first stamp a throwaway placeholder to force a real\n    # edit, then write the real target text. Both assignments exclude\n    # the trailing paragraph mark.\n    $r = $p.Range\n    $r.End = $r.End - 1\n    $r.Text = \".\"\n\n    $r2 = $p.Range\n    $r2.End = $r2.End - 1\n    $r2.Text = $target\n  }\n}\n"}
